$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": two source files were re-generated by the handoff run.
#   0907d6fd-6ff9-467c-bf2b-7c4a9cf7dca3.md -> 06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md
#   e3ec1fca-15d8-45a0-b879-70efba0eba94.md -> ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md
# Status text and handoff datetime also changed.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-47-12 12:47:44"

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-47-12 12:47:44"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md", "", "", "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md", "", "", "ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": source files renamed, new handoff/handback timestamps, and
# the (now-irrelevant) "Dependency From" / duplicate handoff columns F & G
# are cleared out entirely (no more dependency chain between the two files).
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-12 12:47:41"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"

$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-12 12:47:41"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"

$wsZh.Range("F2:G3").Clear()

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md", "", "", "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a06cbc13b215711097bf57b4e1d878cd9a44c76f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.zh-cn.xlf", "", "", "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md", "", "", "ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a06cbc13b215711097bf57b4e1d878cd9a44c76f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.zh-cn.xlf", "", "", "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de": same shape of change as zh-cn, different locale file + times.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-12 12:47:44"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"

$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-12 12:47:44"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"

$wsDe.Range("F2:G3").Clear()

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md", "", "", "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3790e7581cb37d2791d23b0a6601f05602e6e662/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.de-de.xlf", "", "", "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md", "", "", "ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/71bc6dd27e573b6390b52b8e5cbbfa3a765ddc9e/e2e/ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3790e7581cb37d2791d23b0a6601f05602e6e662/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.de-de.xlf", "", "", "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.de-de.xlf")

Write-Output "Localization status report regenerated for handoff."
